$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '28.026.08'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +3.30%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.571.96'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +0.32%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.996'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -1.41%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '212.09'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -1.18%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '23.16'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +5.44%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.250'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +0.35%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0597'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.0880'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +1.59%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.796.79'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '1.571.64'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +1.20%  '

$ws.Range("E14").Value = '  -0.57%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.520'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '27.984.11'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +3.36%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '63.31'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +1.85%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '227.71'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +5.66%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.0₃0704'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '7.44'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +0.95%  '

$ws.Range("E21").Value = '  -1.27%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '4.11'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.54%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '9.31'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +1.25%  '

$ws.Range("E24").Value = '  -0.54%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '151.47'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -1.77%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '15.18'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.91%  '

$ws.Range("E27").Value = '  -0.77%  '

$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '0.996'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.24%  '

$ws.Range("E30").Value = '  -0.16%  '

$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("E32").Value = '  -0.54%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '3.13'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -1.81%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '1.413.42'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -2.05%  '

$ws.Range("E35").Value = '  -1.59%  '

$ws.Range("E36").Value = '  -3.78%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '2.31'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -2.01%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.0167'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.539'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +1.49%  '

$ws.Range("E40").Value = '  +3.18%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.805'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -0.16%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.995'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -1.45%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '5.62'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -3.18%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '1.83'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +5.05%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.969'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -3.20%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '63.45'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -1.79%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '1.708.42'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +0.42%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '86.72'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +1.11%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +3.03%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.0525'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +1.05%  '

$ws.Range("E51").Value = '  -1.84%  '
